$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 372.36667
$ws.Range("J17").Value = 306.67856
$ws.Range("L17").Value = 920.03568
$ws.Range("N17").Value = -1256.03568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1024.48
$ws.Range("I28").Value = 834.4211
$ws.Range("K28").Value = 834.4211
$ws.Range("M28").Value = -349.4211

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1076.4166
$ws.Range("I39").Value = 1190.8889
$ws.Range("K39").Value = 3572.6667
$ws.Range("M39").Value = -3276.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3946.3333
$ws.Range("I43").Value = 4320
$ws.Range("J43").Value = 3829.5625
$ws.Range("K43").Value = 4320
$ws.Range("L43").Value = 3829.5625
$ws.Range("M43").Value = -4251
$ws.Range("N43").Value = -3967.5625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6250
$ws.Range("I86").Value = 7000
$ws.Range("K86").Value = 7000
$ws.Range("M86").Value = -5877

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 6250
$ws.Range("I89").Value = 7000
$ws.Range("K89").Value = 35000
$ws.Range("M89").Value = -29384

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4232.724
$ws.Range("I106").Value = 4583.2856
$ws.Range("J106").Value = 3312.5
$ws.Range("K106").Value = 4583.2856
$ws.Range("L106").Value = 3312.5
$ws.Range("M106").Value = -3952.2856
$ws.Range("N106").Value = -4574.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1333.9375
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1356.2
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 4068.6
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -6284.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3224.8333
$ws.Range("J138").Value = 3593.862
$ws.Range("L138").Value = 10781.586
$ws.Range("N138").Value = -21061.586

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 9787.308000000001
$ws.Range("I141").Value = 10698.125
$ws.Range("J141").Value = 8330
$ws.Range("K141").Value = 32094.375
$ws.Range("L141").Value = 24990
$ws.Range("M141").Value = -26914.375
$ws.Range("N141").Value = -35350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1774.8214
$ws.Range("I32").Value = 1758.5
$ws.Range("J32").Value = 1987
$ws.Range("K32").Value = 1758.5
$ws.Range("L32").Value = 1987
$ws.Range("M32").Value = -1471.5
$ws.Range("N32").Value = -2561

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6266.5454
$ws.Range("I74").Value = 5693.2
$ws.Range("J74").Value = 12000
$ws.Range("K74").Value = 5693.2
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = -4819.2
$ws.Range("N74").Value = -13748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6266.5454
$ws.Range("I77").Value = 5693.2
$ws.Range("J77").Value = 12000
$ws.Range("K77").Value = 28466
$ws.Range("L77").Value = 60000
$ws.Range("M77").Value = -24098
$ws.Range("N77").Value = -68736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1758.3636
$ws.Range("I132").Value = 1758.3636
$ws.Range("K132").Value = 5275.0908
$ws.Range("M132").Value = -2745.0908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = -360
$ws.Range("M8").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7885.0625
$ws.Range("I20").Value = 5205.8184
$ws.Range("J20").Value = 13779.4
$ws.Range("K20").Value = 5205.8184
$ws.Range("L20").Value = 13779.4
$ws.Range("M20").Value = -4958.8184
$ws.Range("N20").Value = -14273.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 9374.125
$ws.Range("I94").Value = 8998.6
$ws.Range("K94").Value = 8998.6
$ws.Range("M94").Value = -8547.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 83355
$ws.Range("I99").Value = 45878.332
$ws.Range("K99").Value = 45878.332
$ws.Range("M99").Value = -44380.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 9340.904
$ws.Range("I107").Value = 8796.117
$ws.Range("J107").Value = 11656.25
$ws.Range("K107").Value = 8796.117
$ws.Range("L107").Value = 11656.25
$ws.Range("M107").Value = -6876.117
$ws.Range("N107").Value = -15496.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5914.2856
$ws.Range("I16").Value = 3975
$ws.Range("J16").Value = 8500
$ws.Range("K16").Value = 3975
$ws.Range("L16").Value = 8500
$ws.Range("M16").Value = -3688
$ws.Range("N16").Value = -9074

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1540
$ws.Range("J94").Value = 1869.4286
$ws.Range("L94").Value = 1869.4286
$ws.Range("N94").Value = -2771.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 5914.2856
$ws.Range("I113").Value = 3975
$ws.Range("J113").Value = 8500
$ws.Range("K113").Value = 3975
$ws.Range("L113").Value = 8500
$ws.Range("M113").Value = -1805
$ws.Range("N113").Value = -12840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4208
$ws.Range("I122").Value = 3727.75
$ws.Range("J122").Value = 4448.125
$ws.Range("K122").Value = 11183.25
$ws.Range("L122").Value = 13344.375
$ws.Range("M122").Value = -8733.25
$ws.Range("N122").Value = -18244.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3852.1667
$ws.Range("I132").Value = 3366.3333
$ws.Range("K132").Value = 10098.9999
$ws.Range("M132").Value = -7568.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 305.3
$ws.Range("I14").Value = 305.3
$ws.Range("K14").Value = 915.9000000000001
$ws.Range("M14").Value = -742.9000000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1055068.9
$ws.Range("I131").Value = 2000985.8
$ws.Range("J131").Value = 4050
$ws.Range("K131").Value = 6002957.4
$ws.Range("L131").Value = 12150
$ws.Range("M131").Value = -5997917.4
$ws.Range("N131").Value = -22230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4401.615
$ws.Range("I122").Value = 3597.15
$ws.Range("K122").Value = 10791.45
$ws.Range("M122").Value = -8341.450000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7740.676
$ws.Range("I132").Value = 6857.069
$ws.Range("J132").Value = 10943.75
$ws.Range("K132").Value = 20571.207
$ws.Range("L132").Value = 32831.25
$ws.Range("M132").Value = -18041.207
$ws.Range("N132").Value = -37891.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2940.6428
$ws.Range("I132").Value = 2474.5386
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 7423.6158
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -4893.6158
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 840.125
$ws.Range("I107").Value = 845.8570999999999
$ws.Range("K107").Value = 2537.5713
$ws.Range("M107").Value = -617.5712999999996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4969.8716
$ws.Range("I132").Value = 3563.5356
$ws.Range("J132").Value = 8549.637000000001
$ws.Range("K132").Value = 10690.6068
$ws.Range("L132").Value = 25648.911
$ws.Range("M132").Value = -8160.606800000001
$ws.Range("N132").Value = -30708.911
